$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-11 with new values (in place) and add new rows 12-21

# Row 2
$ws.Cells.Item(2, 1).Value = "sentence-transformers/all-MiniLM-L6-v2"
$ws.Cells.Item(2, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(2, 3).Value = "{'name': 'sentence-transformers/all-MiniLM-L6-v2', 'batch_size': 100, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:'}"
$ws.Cells.Item(2, 4).Value = 0.5230769230769231
$ws.Cells.Item(2, 5).Value = 0.5692307692307692
$ws.Cells.Item(2, 6).Value = 0.6153846153846154
$ws.Cells.Item(2, 7).Value = 0.676923076923077
$ws.Cells.Item(2, 8).Value = 0.6923076923076923
$ws.Cells.Item(2, 9).Value = 0.5800000000000001

# Row 3
$ws.Cells.Item(3, 1).Value = "mixedbread-ai/mxbai-embed-large-v1"
$ws.Cells.Item(3, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(3, 3).Value = "{'name': 'mixedbread-ai/mxbai-embed-large-v1', 'batch_size': 100, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:'}"
$ws.Cells.Item(3, 4).Value = 0.4307692307692308
$ws.Cells.Item(3, 5).Value = 0.5384615384615384
$ws.Cells.Item(3, 6).Value = 0.6461538461538462
$ws.Cells.Item(3, 7).Value = 0.676923076923077
$ws.Cells.Item(3, 8).Value = 0.6923076923076923
$ws.Cells.Item(3, 9).Value = 0.5312820512820513

# Row 4
$ws.Cells.Item(4, 1).Value = "nvidia/NV-Embed-v2"
$ws.Cells.Item(4, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(4, 3).Value = "{'name': 'nvidia/NV-Embed-v2', 'batch_size': 5, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:', 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True, 'max_length': 32768}}"
$ws.Cells.Item(4, 4).Value = 0.9384615384615385
$ws.Cells.Item(4, 5).Value = 0.9692307692307692
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 0.9641025641025641

# Row 5
$ws.Cells.Item(5, 1).Value = "dunzhang/stella_en_1.5B_v5"
$ws.Cells.Item(5, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(5, 3).Value = "{'name': 'dunzhang/stella_en_1.5B_v5', 'batch_size': 20, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:', 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True}}"
$ws.Cells.Item(5, 4).Value = 0.2769230769230769
$ws.Cells.Item(5, 5).Value = 0.6
$ws.Cells.Item(5, 6).Value = 0.6307692307692307
$ws.Cells.Item(5, 7).Value = 0.7076923076923077
$ws.Cells.Item(5, 8).Value = 0.7076923076923077
$ws.Cells.Item(5, 9).Value = 0.467948717948718

# Row 6
$ws.Cells.Item(6, 1).Value = "amazon.titan-embed-text-v2:0"
$ws.Cells.Item(6, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(6, 3).Value = "{'name': 'amazon.titan-embed-text-v2:0', 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:', 'model_kwargs': {'aws': True, 'aws_creds_file': '/home/ubuntu/Multi-Agent-LLM-System-with-LangGraph-RAG-and-LangChain/config/config.ini', 'aws_config_name': 'BedRock_LLM_API'}}"
$ws.Cells.Item(6, 4).Value = 0.8
$ws.Cells.Item(6, 5).Value = 0.8923076923076924
$ws.Cells.Item(6, 6).Value = 0.9384615384615385
$ws.Cells.Item(6, 7).Value = 0.9384615384615385
$ws.Cells.Item(6, 8).Value = 0.9538461538461539
$ws.Cells.Item(6, 9).Value = 0.8646153846153847

# Row 7
$ws.Cells.Item(7, 1).Value = "sentence-transformers/all-MiniLM-L6-v2"
$ws.Cells.Item(7, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(7, 3).Value = "{'name': 'sentence-transformers/all-MiniLM-L6-v2', 'batch_size': 100}"
$ws.Cells.Item(7, 4).Value = 0.6923076923076923
$ws.Cells.Item(7, 5).Value = 0.7384615384615385
$ws.Cells.Item(7, 6).Value = 0.7384615384615385
$ws.Cells.Item(7, 7).Value = 0.7538461538461538
$ws.Cells.Item(7, 8).Value = 0.7846153846153846
$ws.Cells.Item(7, 9).Value = 0.7253846153846154

# Row 8
$ws.Cells.Item(8, 1).Value = "mixedbread-ai/mxbai-embed-large-v1"
$ws.Cells.Item(8, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(8, 3).Value = "{'name': 'mixedbread-ai/mxbai-embed-large-v1', 'batch_size': 100}"
$ws.Cells.Item(8, 4).Value = 0.8769230769230769
$ws.Cells.Item(8, 5).Value = 0.9692307692307692
$ws.Cells.Item(8, 6).Value = 0.9846153846153847
$ws.Cells.Item(8, 7).Value = 0.9846153846153847
$ws.Cells.Item(8, 8).Value = 0.9846153846153847
$ws.Cells.Item(8, 9).Value = 0.9282051282051281

# Row 9
$ws.Cells.Item(9, 1).Value = "nvidia/NV-Embed-v2"
$ws.Cells.Item(9, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(9, 3).Value = "{'name': 'nvidia/NV-Embed-v2', 'batch_size': 5, 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True, 'max_length': 32768}}"
$ws.Cells.Item(9, 4).Value = 0.8769230769230769
$ws.Cells.Item(9, 5).Value = 0.9384615384615385
$ws.Cells.Item(9, 6).Value = 0.9538461538461539
$ws.Cells.Item(9, 7).Value = 0.9538461538461539
$ws.Cells.Item(9, 8).Value = 0.9538461538461539
$ws.Cells.Item(9, 9).Value = 0.9128205128205129

# Row 10
$ws.Cells.Item(10, 1).Value = "dunzhang/stella_en_1.5B_v5"
$ws.Cells.Item(10, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(10, 3).Value = "{'name': 'dunzhang/stella_en_1.5B_v5', 'batch_size': 20, 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True}}"
$ws.Cells.Item(10, 4).Value = 0.2
$ws.Cells.Item(10, 5).Value = 0.4153846153846154
$ws.Cells.Item(10, 6).Value = 0.5230769230769231
$ws.Cells.Item(10, 7).Value = 0.5846153846153846
$ws.Cells.Item(10, 8).Value = 0.6153846153846154
$ws.Cells.Item(10, 9).Value = 0.3651282051282051

# Row 11
$ws.Cells.Item(11, 1).Value = "amazon.titan-embed-text-v2:0"
$ws.Cells.Item(11, 2).Value = "HuggingFace QA Dataset"
$ws.Cells.Item(11, 3).Value = "{'name': 'amazon.titan-embed-text-v2:0', 'model_kwargs': {'aws': True, 'aws_creds_file': '/home/ubuntu/Multi-Agent-LLM-System-with-LangGraph-RAG-and-LangChain/config/config.ini', 'aws_config_name': 'BedRock_LLM_API'}}"
$ws.Cells.Item(11, 4).Value = 0.8461538461538461
$ws.Cells.Item(11, 5).Value = 0.9076923076923077
$ws.Cells.Item(11, 6).Value = 0.9230769230769231
$ws.Cells.Item(11, 7).Value = 0.9538461538461539
$ws.Cells.Item(11, 8).Value = 0.9538461538461539
$ws.Cells.Item(11, 9).Value = 0.8897435897435897

# Row 12
$ws.Cells.Item(12, 1).Value = "sentence-transformers/all-MiniLM-L6-v2"
$ws.Cells.Item(12, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(12, 3).Value = "{'name': 'sentence-transformers/all-MiniLM-L6-v2', 'batch_size': 100, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:'}"
$ws.Cells.Item(12, 4).Value = 0.6153846153846154
$ws.Cells.Item(12, 5).Value = 0.6923076923076923
$ws.Cells.Item(12, 6).Value = 0.6923076923076923
$ws.Cells.Item(12, 7).Value = 0.6923076923076923
$ws.Cells.Item(12, 8).Value = 0.6923076923076923
$ws.Cells.Item(12, 9).Value = 0.6538461538461539

# Row 13
$ws.Cells.Item(13, 1).Value = "mixedbread-ai/mxbai-embed-large-v1"
$ws.Cells.Item(13, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(13, 3).Value = "{'name': 'mixedbread-ai/mxbai-embed-large-v1', 'batch_size': 100, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:'}"
$ws.Cells.Item(13, 4).Value = 0.6923076923076923
$ws.Cells.Item(13, 5).Value = 0.8461538461538461
$ws.Cells.Item(13, 6).Value = 0.9230769230769231
$ws.Cells.Item(13, 7).Value = 0.9230769230769231
$ws.Cells.Item(13, 8).Value = 0.9230769230769231
$ws.Cells.Item(13, 9).Value = 0.7948717948717949

# Row 14
$ws.Cells.Item(14, 1).Value = "nvidia/NV-Embed-v2"
$ws.Cells.Item(14, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(14, 3).Value = "{'name': 'nvidia/NV-Embed-v2', 'batch_size': 5, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:', 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True, 'max_length': 32768}}"
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(14, 9).Value = 1

# Row 15
$ws.Cells.Item(15, 1).Value = "dunzhang/stella_en_1.5B_v5"
$ws.Cells.Item(15, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(15, 3).Value = "{'name': 'dunzhang/stella_en_1.5B_v5', 'batch_size': 20, 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:', 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True}}"
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 1

# Row 16
$ws.Cells.Item(16, 1).Value = "amazon.titan-embed-text-v2:0"
$ws.Cells.Item(16, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(16, 3).Value = "{'name': 'amazon.titan-embed-text-v2:0', 'instruction': 'Instruct: Represent this passage for retrieval in response to relevant questions.\nQuery:', 'query_instruction': 'Instruct: Given a query, find the most relevant passages that can provide the answer.\nPassage:', 'model_kwargs': {'aws': True, 'aws_creds_file': '/home/ubuntu/Multi-Agent-LLM-System-with-LangGraph-RAG-and-LangChain/config/config.ini', 'aws_config_name': 'BedRock_LLM_API'}}"
$ws.Cells.Item(16, 4).Value = 0.9230769230769231
$ws.Cells.Item(16, 5).Value = 0.9230769230769231
$ws.Cells.Item(16, 6).Value = 0.9230769230769231
$ws.Cells.Item(16, 7).Value = 0.9230769230769231
$ws.Cells.Item(16, 8).Value = 0.9230769230769231
$ws.Cells.Item(16, 9).Value = 0.9230769230769231

# Row 17
$ws.Cells.Item(17, 1).Value = "sentence-transformers/all-MiniLM-L6-v2"
$ws.Cells.Item(17, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(17, 3).Value = "{'name': 'sentence-transformers/all-MiniLM-L6-v2', 'batch_size': 100}"
$ws.Cells.Item(17, 4).Value = 0.6923076923076923
$ws.Cells.Item(17, 5).Value = 0.7692307692307693
$ws.Cells.Item(17, 6).Value = 0.8461538461538461
$ws.Cells.Item(17, 7).Value = 0.8461538461538461
$ws.Cells.Item(17, 8).Value = 0.8461538461538461
$ws.Cells.Item(17, 9).Value = 0.7564102564102563

# Row 18
$ws.Cells.Item(18, 1).Value = "mixedbread-ai/mxbai-embed-large-v1"
$ws.Cells.Item(18, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(18, 3).Value = "{'name': 'mixedbread-ai/mxbai-embed-large-v1', 'batch_size': 100}"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 1

# Row 19
$ws.Cells.Item(19, 1).Value = "nvidia/NV-Embed-v2"
$ws.Cells.Item(19, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(19, 3).Value = "{'name': 'nvidia/NV-Embed-v2', 'batch_size': 5, 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True, 'max_length': 32768}}"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 1
$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(19, 9).Value = 1

# Row 20
$ws.Cells.Item(20, 1).Value = "dunzhang/stella_en_1.5B_v5"
$ws.Cells.Item(20, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(20, 3).Value = "{'name': 'dunzhang/stella_en_1.5B_v5', 'batch_size': 20, 'model_kwargs': {'trust_remote_code': True, 'load_in_8bit': True}}"
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(20, 9).Value = 1

# Row 21
$ws.Cells.Item(21, 1).Value = "amazon.titan-embed-text-v2:0"
$ws.Cells.Item(21, 2).Value = "PubMed filtered Dataset"
$ws.Cells.Item(21, 3).Value = "{'name': 'amazon.titan-embed-text-v2:0', 'model_kwargs': {'aws': True, 'aws_creds_file': '/home/ubuntu/Multi-Agent-LLM-System-with-LangGraph-RAG-and-LangChain/config/config.ini', 'aws_config_name': 'BedRock_LLM_API'}}"
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 1
$ws.Cells.Item(21, 9).Value = 1
